$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Cells.Item(2, 12).Value = 2026  # L2 was 2006
$ws.Cells.Item(3, 12).Value = 2060  # L3 was 2043
$ws.Cells.Item(4, 8).Value = 1754  # H4 was 1752
$ws.Cells.Item(4, 10).Value = 1864  # J4 was 1862
$ws.Cells.Item(4, 11).Value = 1760  # K4 was 1759
$ws.Cells.Item(4, 12).Value = 572  # L4 was 567
$ws.Cells.Item(6, 11).Value = 9124  # K6 was 9125
$ws.Cells.Item(6, 12).Value = 1843  # L6 was 1830
$ws.Cells.Item(7, 8).Value = 26067  # H7 was 26065
$ws.Cells.Item(7, 10).Value = 29337  # J7 was 29335
$ws.Cells.Item(7, 12).Value = 6618  # L7 was 6563

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Cells.Item(6, 12).Value = 28  # L6 was 27
$ws.Cells.Item(7, 12).Value = 79  # L7 was 78

$ws = $wb.Worksheets.Item('Austin')
$ws.Cells.Item(2, 12).Value = 117  # L2 was 116
$ws.Cells.Item(6, 12).Value = 111  # L6 was 110
$ws.Cells.Item(7, 12).Value = 416  # L7 was 414

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Cells.Item(2, 12).Value = 50  # L2 was 49
$ws.Cells.Item(3, 12).Value = 64  # L3 was 63
$ws.Cells.Item(7, 12).Value = 156  # L7 was 154

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Cells.Item(3, 12).Value = 100  # L3 was 99
$ws.Cells.Item(6, 12).Value = 101  # L6 was 100
$ws.Cells.Item(7, 12).Value = 294  # L7 was 292

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Cells.Item(6, 12).Value = 20  # L6 was 19
$ws.Cells.Item(7, 12).Value = 98  # L7 was 97

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Cells.Item(3, 12).Value = 67  # L3 was 66
$ws.Cells.Item(7, 12).Value = 237  # L7 was 236

$ws = $wb.Worksheets.Item('New City')
$ws.Cells.Item(2, 12).Value = 48  # L2 was 45
$ws.Cells.Item(3, 12).Value = 40  # L3 was 38
$ws.Cells.Item(7, 12).Value = 129  # L7 was 124

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Cells.Item(6, 12).Value = 21  # L6 was 19
$ws.Cells.Item(7, 12).Value = 102  # L7 was 100

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Cells.Item(2, 12).Value = 50  # L2 was 48
$ws.Cells.Item(7, 12).Value = 222  # L7 was 221
$ws.Cells.Item(8, 12).Value = 416  # L8 was 414
$ws.Cells.Item(9, 12).Value = 40  # L9 was 38
$ws.Cells.Item(10, 12).Value = 45  # L10 was 43
$ws.Cells.Item(18, 12).Value = 48  # L18 was 46
$ws.Cells.Item(19, 12).Value = 187  # L19 was 185
$ws.Cells.Item(29, 12).Value = 340  # L29 was 338
$ws.Cells.Item(33, 12).Value = 294  # L33 was 292
$ws.Cells.Item(34, 12).Value = 42  # L34 was 41
$ws.Cells.Item(37, 12).Value = 237  # L37 was 236
$ws.Cells.Item(42, 12).Value = 206  # L42 was 205
$ws.Cells.Item(47, 12).Value = 45  # L47 was 44
$ws.Cells.Item(51, 12).Value = 78  # L51 was 77
$ws.Cells.Item(52, 12).Value = 135  # L52 was 133
$ws.Cells.Item(53, 12).Value = 79  # L53 was 78
$ws.Cells.Item(54, 12).Value = 139  # L54 was 138
$ws.Cells.Item(63, 8).Value = 303  # H63 was 301
$ws.Cells.Item(63, 10).Value = 216  # J63 was 214
$ws.Cells.Item(63, 11).Value = 91  # K63 was 90
$ws.Cells.Item(63, 12).Value = 26  # L63 was 22
$ws.Cells.Item(64, 12).Value = 49  # L64 was 48
$ws.Cells.Item(65, 12).Value = 129  # L65 was 124
$ws.Cells.Item(67, 12).Value = 236  # L67 was 235
$ws.Cells.Item(76, 12).Value = 69  # L76 was 68
$ws.Cells.Item(78, 12).Value = 93  # L78 was 91
$ws.Cells.Item(83, 12).Value = 156  # L83 was 154
$ws.Cells.Item(84, 12).Value = 69  # L84 was 68
$ws.Cells.Item(85, 12).Value = 349  # L85 was 348
$ws.Cells.Item(86, 12).Value = 51  # L86 was 50
$ws.Cells.Item(87, 12).Value = 19  # L87 was 18
$ws.Cells.Item(89, 12).Value = 85  # L89 was 84
$ws.Cells.Item(90, 11).Value = 261  # K90 was 262
$ws.Cells.Item(93, 12).Value = 36  # L93 was 35
$ws.Cells.Item(94, 12).Value = 79  # L94 was 78
$ws.Cells.Item(95, 12).Value = 98  # L95 was 97
$ws.Cells.Item(96, 12).Value = 60  # L96 was 59
$ws.Cells.Item(98, 12).Value = 49  # L98 was 47
$ws.Cells.Item(99, 12).Value = 102  # L99 was 100
$ws.Cells.Item(101, 8).Value = 26067  # H101 was 26065
$ws.Cells.Item(101, 10).Value = 29337  # J101 was 29335
$ws.Cells.Item(101, 12).Value = 6618  # L101 was 6563

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Cells.Item(2, 12).Value = 67  # L2 was 66
$ws.Cells.Item(7, 12).Value = 236  # L7 was 235

$ws = $wb.Worksheets.Item('South Deering')
$ws.Cells.Item(6, 12).Value = 15  # L6 was 14
$ws.Cells.Item(7, 12).Value = 69  # L7 was 68

$ws = $wb.Worksheets.Item('Loop')
$ws.Cells.Item(3, 12).Value = 26  # L3 was 25
$ws.Cells.Item(7, 12).Value = 139  # L7 was 138

$ws = $wb.Worksheets.Item('Englewood')
$ws.Cells.Item(3, 12).Value = 122  # L3 was 121
$ws.Cells.Item(6, 12).Value = 90  # L6 was 89
$ws.Cells.Item(7, 12).Value = 340  # L7 was 338

$ws = $wb.Worksheets.Item('Chatham')
$ws.Cells.Item(3, 12).Value = 58  # L3 was 57
$ws.Cells.Item(6, 12).Value = 59  # L6 was 58
$ws.Cells.Item(7, 12).Value = 187  # L7 was 185

$ws = $wb.Worksheets.Item('River North')
$ws.Cells.Item(2, 12).Value = 12  # L2 was 11
$ws.Cells.Item(7, 12).Value = 69  # L7 was 68

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Cells.Item(2, 12).Value = 54  # L2 was 53
$ws.Cells.Item(7, 12).Value = 206  # L7 was 205

$ws = $wb.Worksheets.Item('Avondale')
$ws.Cells.Item(2, 12).Value = 20  # L2 was 18
$ws.Cells.Item(7, 12).Value = 45  # L7 was 43

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Cells.Item(4, 12).Value = 12  # L4 was 11
$ws.Cells.Item(6, 12).Value = 26  # L6 was 25
$ws.Cells.Item(7, 12).Value = 93  # L7 was 91

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Cells.Item(2, 12).Value = 26  # L2 was 25
$ws.Cells.Item(7, 12).Value = 60  # L7 was 59

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Cells.Item(2, 12).Value = 17  # L2 was 16
$ws.Cells.Item(7, 12).Value = 49  # L7 was 48

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Cells.Item(3, 12).Value = 18  # L3 was 17
$ws.Cells.Item(6, 12).Value = 7  # L6 was 6
$ws.Cells.Item(7, 12).Value = 48  # L7 was 46

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Cells.Item(2, 12).Value = 12  # L2 was 11
$ws.Cells.Item(7, 12).Value = 36  # L7 was 35

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Cells.Item(3, 12).Value = 70  # L3 was 69
$ws.Cells.Item(7, 12).Value = 222  # L7 was 221

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Cells.Item(6, 12).Value = 16  # L6 was 15
$ws.Cells.Item(7, 12).Value = 42  # L7 was 41

$ws = $wb.Worksheets.Item('West Loop')
$ws.Cells.Item(3, 12).Value = 20  # L3 was 19
$ws.Cells.Item(7, 12).Value = 79  # L7 was 78

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Cells.Item(3, 12).Value = 16  # L3 was 15
$ws.Cells.Item(7, 12).Value = 45  # L7 was 44

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Cells.Item(3, 12).Value = 9  # L3 was 8
$ws.Cells.Item(4, 12).Value = 4  # L4 was 3
$ws.Cells.Item(7, 12).Value = 49  # L7 was 47

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Cells.Item(2, 12).Value = 10  # L2 was 9
$ws.Cells.Item(3, 12).Value = 18  # L3 was 17
$ws.Cells.Item(7, 12).Value = 40  # L7 was 38

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Cells.Item(2, 12).Value = 14  # L2 was 12
$ws.Cells.Item(7, 12).Value = 50  # L7 was 48

$ws = $wb.Worksheets.Item('Uptown')
$ws.Cells.Item(2, 12).Value = 32  # L2 was 31
$ws.Cells.Item(7, 12).Value = 85  # L7 was 84

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Cells.Item(2, 12).Value = 9  # L2 was 8
$ws.Cells.Item(7, 12).Value = 51  # L7 was 50

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Cells.Item(6, 11).Value = 69  # K6 was 70
$ws.Cells.Item(7, 11).Value = 261  # K7 was 262

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Cells.Item(3, 12).Value = 24  # L3 was 23
$ws.Cells.Item(7, 12).Value = 78  # L7 was 77

$ws = $wb.Worksheets.Item('South Shore')
$ws.Cells.Item(2, 12).Value = 108  # L2 was 107
$ws.Cells.Item(7, 12).Value = 349  # L7 was 348

$ws = $wb.Worksheets.Item('Little Village')
$ws.Cells.Item(3, 12).Value = 41  # L3 was 39
$ws.Cells.Item(7, 12).Value = 135  # L7 was 133

$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Cells.Item(2, 12).Value = 6  # L2 was 5
$ws.Cells.Item(7, 12).Value = 19  # L7 was 18
